$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit swaps the content of row 2 and row 3 for the columns that
# differ between the two records (A, B, D, E, F, G, H, P, Q, R, Y, AA).
# Every other column already holds identical values in both rows, so a
# full swap of just these columns reproduces the diff.

function Set-DateTextValue {
    param($range, $text)
    # Writing a date-like string straight into Value2 makes Excel parse it
    # as a real date (and apply a date number format). Force a text
    # number format first so the literal string is preserved, then revert
    # the format to the default "Normal" style so no stray style lingers
    # on the cell.
    $range.NumberFormat = "@"
    $range.Value2 = $text
    $range.Style = "Normal"
}

$cols = @("A","B","D","E","F","G","H","P","Q","R")
$dateCols = @("Y","AA")

$row2 = @{}
$row3 = @{}
foreach ($col in $cols) {
    $row2[$col] = $ws.Range("$col`2").Value2
    $row3[$col] = $ws.Range("$col`3").Value2
}
foreach ($col in $dateCols) {
    $row2[$col] = $ws.Range("$col`2").Text
    $row3[$col] = $ws.Range("$col`3").Text
}

# Write row 3's original values into row 2, and row 2's original values into row 3.
foreach ($col in $cols) {
    $ws.Range("$col`2").Value2 = $row3[$col]
    $ws.Range("$col`3").Value2 = $row2[$col]
}
foreach ($col in $dateCols) {
    Set-DateTextValue $ws.Range("$col`2") $row3[$col]
    Set-DateTextValue $ws.Range("$col`3") $row2[$col]
}
